$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 281, pushing existing rows 281:353 down to 284:356
$ws.Rows("281:283").Insert()

# Populate the newly inserted rows with the new "September Snow" variety data
# Row 281 - Especial
$ws.Cells.Item(281, 1).Value = 11
$ws.Cells.Item(281, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(281, 3).Value = "Bíobío"
$ws.Cells.Item(281, 4).Value = 44988
$ws.Cells.Item(281, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(281, 5).Value = 8
$ws.Cells.Item(281, 6).Value = "Fruta"
$ws.Cells.Item(281, 7).Value = 100103
$ws.Cells.Item(281, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(281, 9).Value = 100103004
$ws.Cells.Item(281, 10).Value = "Durazno"
$ws.Cells.Item(281, 11).Value = "September Snow"
$ws.Cells.Item(281, 12).Value = "Especial"
$ws.Cells.Item(281, 13).Value = 50
$ws.Cells.Item(281, 14).Value = 16000
$ws.Cells.Item(281, 15).Value = 16000
$ws.Cells.Item(281, 16).Value = 16000
$ws.Cells.Item(281, 17).Value = "`$/caja 15 kilos empedrada"
$ws.Cells.Item(281, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(281, 19).Value = 1067
$ws.Cells.Item(281, 20).Value = 15

# Row 282 - Primera
$ws.Cells.Item(282, 1).Value = 11
$ws.Cells.Item(282, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(282, 3).Value = "Bíobío"
$ws.Cells.Item(282, 4).Value = 44988
$ws.Cells.Item(282, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(282, 5).Value = 8
$ws.Cells.Item(282, 6).Value = "Fruta"
$ws.Cells.Item(282, 7).Value = 100103
$ws.Cells.Item(282, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(282, 9).Value = 100103004
$ws.Cells.Item(282, 10).Value = "Durazno"
$ws.Cells.Item(282, 11).Value = "September Snow"
$ws.Cells.Item(282, 12).Value = "Primera"
$ws.Cells.Item(282, 13).Value = 50
$ws.Cells.Item(282, 14).Value = 14000
$ws.Cells.Item(282, 15).Value = 14000
$ws.Cells.Item(282, 16).Value = 14000
$ws.Cells.Item(282, 17).Value = "`$/caja 15 kilos empedrada"
$ws.Cells.Item(282, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(282, 19).Value = 933
$ws.Cells.Item(282, 20).Value = 15

# Row 283 - Segunda
$ws.Cells.Item(283, 1).Value = 11
$ws.Cells.Item(283, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(283, 3).Value = "Bíobío"
$ws.Cells.Item(283, 4).Value = 44988
$ws.Cells.Item(283, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(283, 5).Value = 8
$ws.Cells.Item(283, 6).Value = "Fruta"
$ws.Cells.Item(283, 7).Value = 100103
$ws.Cells.Item(283, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(283, 9).Value = 100103004
$ws.Cells.Item(283, 10).Value = "Durazno"
$ws.Cells.Item(283, 11).Value = "September Snow"
$ws.Cells.Item(283, 12).Value = "Segunda"
$ws.Cells.Item(283, 13).Value = 50
$ws.Cells.Item(283, 14).Value = 12000
$ws.Cells.Item(283, 15).Value = 12000
$ws.Cells.Item(283, 16).Value = 12000
$ws.Cells.Item(283, 17).Value = "`$/caja 15 kilos empedrada"
$ws.Cells.Item(283, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(283, 19).Value = 800
$ws.Cells.Item(283, 20).Value = 15

Write-Host "Done"
